# Apply "Doing Updates for Financials" edits to the RGS worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RGS")

# Row 14 - Non Recurring: D/E switch from "NA" text to actual numbers, F updated
$ws.Range("D14").Value = 13800
$ws.Range("E14").Value = 16300
$ws.Range("F14").Value = 7100

# Row 15 - Others
$ws.Range("D15").Value = 47100
$ws.Range("E15").Value = 110500
$ws.Range("F15").Value = 113300

# Row 45 - Other Current Assets
$ws.Range("D45").Value = 76200

# Row 48 - Property Plant and Equipment
$ws.Range("D48").Value = 211700

# Row 49 - Goodwill
$ws.Range("D49").Value = 835800

# Row 59 - Other Current Liabilities
$ws.Range("D59").Value = 198300

# Row 60 - Total Current Liabilities
$ws.Range("D60").Value = 158500

# Row 62 - Other Liabilities
$ws.Range("D62").Value = 229700

# Row 66 - Total Liabilities
$ws.Range("D66").Value = 370300

# Row 72 - Retained Earnings
$ws.Range("D72").Value = 280100

# Row 76 - Total Stockholder Equity
$ws.Range("D76").Value = 486400
